$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------
# 1) Insert "aws.ses" into the target list (column A), keeping it
#    alphabetically sorted between "aws.s3" and "base".
# ---------------------------------------------------------------
for ($r = 26; $r -ge 3; $r--) {
    $v = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 1).Value2 = $v
}
$ws.Cells.Item(3, 1).Value2 = "aws.ses"

# ---------------------------------------------------------------
# 2) Insert a new column at C (for the aws.ses commands), shifting
#    all the existing category columns C..Z one place right (D..AA)
# ---------------------------------------------------------------
$ws.Columns.Item(3).Insert()

# ---------------------------------------------------------------
# 3) Populate the new aws.ses column with its header + two commands
# ---------------------------------------------------------------
$ws.Range("C1").Value2 = "aws.ses"
$ws.Range("C2").Value2 = "sendMail(profile,to,subject,body)"
$ws.Range("C3").Value2 = "sendTextMail(profile,to,subject,body)"

# ---------------------------------------------------------------
# 4) Fix up the named ranges that now point to the wrong column
#    (everything that used to live in C..Z) plus "target" (new row)
#    and finally register the brand-new "aws.ses" name.
# ---------------------------------------------------------------
$wb.Names.Item("base").Delete()
$wb.Names.Item("csv").Delete()
$wb.Names.Item("desktop").Delete()
$wb.Names.Item("excel").Delete()
$wb.Names.Item("external").Delete()
$wb.Names.Item("image").Delete()
$wb.Names.Item("io").Delete()
$wb.Names.Item("jms").Delete()
$wb.Names.Item("json").Delete()
$wb.Names.Item("mail").Delete()
$wb.Names.Item("number").Delete()
$wb.Names.Item("pdf").Delete()
$wb.Names.Item("rdbms").Delete()
$wb.Names.Item("redis").Delete()
$wb.Names.Item("sms").Delete()
$wb.Names.Item("sound").Delete()
$wb.Names.Item("ssh").Delete()
$wb.Names.Item("step").Delete()
$wb.Names.Item("target").Delete()
$wb.Names.Item("web").Delete()
$wb.Names.Item("webalert").Delete()
$wb.Names.Item("webcookie").Delete()
$wb.Names.Item("ws").Delete()
$wb.Names.Item("ws.async").Delete()
$wb.Names.Item("xml").Delete()

$wb.Names.Add("base", "='#system'!`$D`$2:`$D`$36")
$wb.Names.Add("csv", "='#system'!`$E`$2:`$E`$5")
$wb.Names.Add("desktop", "='#system'!`$F`$2:`$F`$92")
$wb.Names.Add("excel", "='#system'!`$G`$2:`$G`$14")
$wb.Names.Add("external", "='#system'!`$H`$2:`$H`$3")
$wb.Names.Add("image", "='#system'!`$I`$2:`$I`$5")
$wb.Names.Add("io", "='#system'!`$J`$2:`$J`$24")
$wb.Names.Add("jms", "='#system'!`$K`$2:`$K`$4")
$wb.Names.Add("json", "='#system'!`$L`$2:`$L`$14")
$wb.Names.Add("mail", "='#system'!`$M`$2:`$M`$2")
$wb.Names.Add("number", "='#system'!`$N`$2:`$N`$15")
$wb.Names.Add("pdf", "='#system'!`$O`$2:`$O`$16")
$wb.Names.Add("rdbms", "='#system'!`$P`$2:`$P`$7")
$wb.Names.Add("redis", "='#system'!`$Q`$2:`$Q`$10")
$wb.Names.Add("sms", "='#system'!`$R`$2:`$R`$2")
$wb.Names.Add("sound", "='#system'!`$S`$2:`$S`$5")
$wb.Names.Add("ssh", "='#system'!`$T`$2:`$T`$9")
$wb.Names.Add("step", "='#system'!`$U`$2:`$U`$4")
$wb.Names.Add("target", "='#system'!`$A`$2:`$A`$27")
$wb.Names.Add("web", "='#system'!`$V`$2:`$V`$117")
$wb.Names.Add("webalert", "='#system'!`$W`$2:`$W`$8")
$wb.Names.Add("webcookie", "='#system'!`$X`$2:`$X`$8")
$wb.Names.Add("ws", "='#system'!`$Y`$2:`$Y`$17")
$wb.Names.Add("ws.async", "='#system'!`$Z`$2:`$Z`$8")
$wb.Names.Add("xml", "='#system'!`$AA`$2:`$AA`$11")
$wb.Names.Add("aws.ses", "='#system'!`$C`$2:`$C`$3")
